$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 13, pushing existing row 13 (and below) down.
$ws.Rows.Item(13).Insert()

# Fill in the new data row 13 (continuing the time-log pattern)
$ws.Range("A13").Value = 2014
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 21
$ws.Range("D13").Value = 0.5625
$ws.Range("E13").Value = 0.66666666666666663

# Copy style (number format) from D12/E12/F12/G12 into the new row 13 cells
$ws.Range("D12:G12").Copy()
$ws.Range("D13:G13").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Extend the formula down to F13 (each row keeps its own relative formula)
$ws.Range("F13").Formula = "=(E13-D13)*24*60"

# Update the SUM formula row (now row 15) to include the new row 14 (the blank separator)
$ws.Range("F15").Formula = "=SUM(F2:F14)"

# Set active selection to A14
$ws.Range("A14").Select()
